$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 203, shifting existing rows 203:235 down to 204:236
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with data
$ws.Cells.Item(203, 1).Value = 11
$ws.Cells.Item(203, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value = "Bíobío"
$ws.Cells.Item(203, 4).Value = 45211
$ws.Cells.Item(203, 5).Value = 8
$ws.Cells.Item(203, 6).Value = 100112021
$ws.Cells.Item(203, 7).Value = "Ají"
$ws.Cells.Item(203, 8).Value = "Inferno"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 10
$ws.Cells.Item(203, 11).Value = 32000
$ws.Cells.Item(203, 12).Value = 32000
$ws.Cells.Item(203, 13).Value = 32000
$ws.Cells.Item(203, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(203, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(203, 16).Value = 3200
$ws.Cells.Item(203, 17).Value = 10
$ws.Cells.Item(203, 18).Value = "Hortaliza"

# Ensure the date-cell style/format matches the rest of column D
$ws.Cells.Item(203, 4).NumberFormat = $ws.Cells.Item(204, 4).NumberFormat
